# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -5
$ws.Range("F6").Value = -3
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = -5
$ws.Range("F10").Value = -4
$ws.Range("F12").Value = 3
